$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = 8
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -6
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = -2
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = -2
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = 1
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = -2
